{"js": "// Add 2 \"rapport d'entretien\" entries (commit: \"Ajout de 2 rapports d entretien du 17 mars\")\n//\n// 1) The stray \"_GoBack\" bookmark that sat on the \"Windows account expired\"\n//    paragraph is removed.\n// 2) The \"Mise \u00e0 jour menu\" paragraph gets a new run \"s\" appended (so the\n//    visible text becomes \"Mise \u00e0 jour menus\") and the \"_GoBack\" bookmark is\n//    now placed at the end of that paragraph.\n// 3) A brand-new paragraph \"Menu Eldora en erreur \" is inserted right after it.\n\nconst doc = context.document;\nconst body = doc.body;\n\n// --- Step 1: drop the old \"_GoBack\" bookmark wherever it currently lives ---\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- Step 2: locate the \"Mise \u00e0 jour menu\" paragraph ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Mise \u00e0 jour menu\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Paragraph \"Mise \u00e0 jour menu\" not found');\n}\n\n// --- Step 3: append \"s\" as its own run (w:r) at the end of the paragraph ---\nconst tailRange = target.getRange(Word.RangeLocation.end);\nconst runOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships></pkg:xmlData></pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>s</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part>\n</pkg:package>`;\ntailRange.insertOoxml(runOoxml, Word.InsertLocation.end);\nawait context.sync();\n\n// --- Step 4: re-add the \"_GoBack\" bookmark, now at the end of this paragraph ---\nconst newTailRange = target.getRange(Word.RangeLocation.end);\nnewTailRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- Step 5: insert the new paragraph right after it ---\ntarget.insertParagraph(\"Menu Eldora en erreur \", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Add 2 \"rapport d'entretien\" entries (commit: \"Ajout de 2 rapports d entretien du 17 mars\")\n#\n# 1) The stray \"_GoBack\" bookmark that sat on the \"Windows account expired\"\n#    paragraph is removed.\n# 2) The \"Mise \u00e0 jour menu\" paragraph gets a new run \"s\" appended (so the\n#    visible text becomes \"Mise \u00e0 jour menus\") and the \"_GoBack\" bookmark is\n#    now placed at the end of that paragraph.\n# 3) A brand-new paragraph \"Menu Eldora en erreur \" is inserted right after it.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: drop the old \"_GoBack\" bookmark wherever it currently lives ---\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- Step 2: locate the \"Mise \u00e0 jour menu\" paragraph ---\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\") -eq \"Mise \u00e0 jour menu\") {\n        $target = $p\n        break\n    }\n}\n\n$r = $target.Range\n$contentRange = $d.Range($r.Start, $r.End - 1)   # exclude the paragraph mark\n$insertPos = $contentRange.End\n\n# --- Step 3: append \"s\" as its own run (w:r) at the end of the paragraph ---\n# A scratch paragraph is used so the new text picks up \"clean\" (unstyled)\n# formatting and Word is forced to start a brand-new run instead of just\n# growing the existing \"Mise \u00e0 jour menu\" run.\n$scratchPara = $d.Paragraphs.Add()\n$scratchRange = $scratchPara.Range\n$scratchRange.Text = \"s\"\n$scratchTextRange = $d.Range($scratchRange.Start, $scratchRange.Start + 1)\n$formattedS = $scratchTextRange.FormattedText\n$insertionPoint = $d.Range($insertPos, $insertPos)\n$insertionPoint.FormattedText = $formattedS\n$scratchPara.Range.Delete()\n\n$afterS = $insertPos + 1\n\n# --- Step 4: re-add the \"_GoBack\" bookmark, now at the end of this paragraph ---\n# A temporary placeholder character is inserted first so the bookmark's\n# (collapsed) insertion point is not exactly the paragraph's end-of-text\n# position, then the placeholder is removed once the bookmark is in place.\n$placeholder = $d.Range($afterS, $afterS)\n$placeholder.InsertAfter(\"X\")\n$bmPoint = $d.Range($afterS, $afterS)\n$d.Bookmarks.Add(\"_GoBack\", $bmPoint)\n$placeholderRange = $d.Range($afterS, $afterS + 1)\n$placeholderRange.Delete()\n\n# --- Step 5: insert the new paragraph right after it ---\n$finalContentRange = $d.Range($r.Start, $afterS)\n$finalContentRange.InsertParagraphAfter()\n$newPara = $target.Next()\n$newPara.Range.Text = \"Menu Eldora en erreur \"\n"}
